$d = $word.ActiveDocument

# wdAlignParagraphLeft = 0. Re-asserting "left" alignment on paragraphs that
# are already left-aligned (the implicit default) makes the writer omit the
# <w:jc w:val="left"/> attribute entirely, matching AlignDefault behavior.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Format.Alignment -eq 0) {
        $p.Format.Alignment = 0
    }
}
